# Adds newly-identified 84000/BDRC person matches, pulled from the 84000
# XML, for texts that already had some BDRC data. A new intermediate
# "Sheet1" (the scratch sheet used while cross-checking against Wikidata)
# is inserted between the two existing sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# The 13 newly-identified (84000 name / eft id) pairs, in the order they
# were first typed -- this also controls shared-string allocation order.
# ---------------------------------------------------------------------
$entries = @(
    @("Gö Chödrup", "eft:g-ch-drup"),
    @("wang phab zhwun (wang phan zhun)", "eft:wang-phab-zhwun-wang-phan-zhun-"),
    @("dge ba'i blo gros", "eft:dge-ba-i-blo-gros"),
    @("rgya mtsho'i sde", "eft:rgya-mtsho-i-sde"),
    @("Thönmi Sambhoṭa", "eft:th-nmi-sambhota"),
    @("Tsultrim Gyaltsen", "eft:tsultrim-gyaltsen"),
    @("Shang Buchikpa", "eft:shang-buchikpa"),
    @("Sherap Ö", "eft:sherap-"),
    @("Paṇḍita Dharmākara", "eft:pandita-dharmakara"),
    @("Lotsāwa Zangkyong (bzang skyong)", "eft:lotsawa-zangkyong-bzang-skyong-"),
    @("Nyen Lotsawa Darma Drak", "eft:nyen-lotsawa-darma-drak"),
    @("Patsap Nyima Drak [?]", "eft:patsap-nyima-drak-"),
    @("vajrvisramitra", "eft:vajrvisramitra")
)

# ---------------------------------------------------------------------
# 1) WD_person_matches: append rows 44-56 (columns C/D/E), unresolved
#    BDRC ID marked "?" for every one of them.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("WD_person_matches")

$r = 44
foreach ($pair in $entries) {
    $ws1.Cells.Item($r, 3).Value = $pair[0]
    $ws1.Cells.Item($r, 4).Value = $pair[1]
    $ws1.Cells.Item($r, 5).Value = "?"
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2) Insert a new scratch sheet ("Sheet1") between the two existing
#    sheets, used to cross-check the first four entries against Wikidata.
# ---------------------------------------------------------------------
$new = $wb.Worksheets.Add()
$new.Name = "Sheet1"

$r = 12
foreach ($pair in $entries) {
    $new.Cells.Item($r, 3).Value = $pair[0]
    $new.Cells.Item($r, 4).Value = $pair[1]
    $r = $r + 1
}

# Matching Wikidata person IDs found for the first four rows -- typed in
# this particular order (again controlling shared-string allocation).
$new.Cells.Item(15, 5).Value = "P8277"
$new.Cells.Item(14, 5).Value = "P8278"
$new.Cells.Item(12, 5).Value = "P8221"
$new.Cells.Item(13, 5).Value = "P8276"

# E15 was pasted in from the Wikidata web page, carrying its own font.
$new.Range("E15").Font.Name = "Roboto"
$new.Range("E15").Font.Size = 13
$new.Range("E15").Font.Color = 3421236
$new.Rows.Item(15).RowHeight = 17

$new.Range("E14").Select()

# ---------------------------------------------------------------------
# 3) Tidy up selections on the other two sheets, and leave
#    WD_person_matches as the active tab.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("previously_identified_matches")
$ws3.Range("F92").Select()

$ws1.Activate()
$ws1.Range("E56").Select()
